# The sheet contains one new weekly price record that needs to be inserted
# as row 333 (pushing the former rows 333:405 down to 334:406), growing the
# used range from A1:R405 to A1:R406.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at 333; Excel shifts rows 333:405 -> 334:406
# and extends the sheet dimension automatically.
$ws.Rows("333:333").Insert()

# Populate the newly inserted row 333 with the new record's data.
$ws.Cells.Item(333, 1).Value  = 6
$ws.Cells.Item(333, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(333, 3).Value  = "Metropolitana"
$ws.Cells.Item(333, 4).Value  = 44637
$ws.Cells.Item(333, 5).Value  = 13
$ws.Cells.Item(333, 6).Value  = 100112043
$ws.Cells.Item(333, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(333, 8).Value  = "Sin especificar"
$ws.Cells.Item(333, 9).Value  = "Primera"
$ws.Cells.Item(333, 10).Value = 400
$ws.Cells.Item(333, 11).Value = 17000
$ws.Cells.Item(333, 12).Value = 18000
$ws.Cells.Item(333, 13).Value = 17575
$ws.Cells.Item(333, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(333, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(333, 16).Value = 293
$ws.Cells.Item(333, 17).Value = 60
$ws.Cells.Item(333, 18).Value = "Hortaliza"
